$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Make sure every cell in the table (A1:G18) uses the "text + wrap" style
# that is already used for most of the data columns (numFmt "@", wrapText).
# This upgrades column A (which previously used the plain wrap style) to
# match the rest of the table.
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:G18")
$tableRange.NumberFormat = "@"
$tableRange.WrapText = $true

# ---------------------------------------------------------------------------
# Fill in the author / OOP / actively-developed / latest-version / latest-
# version-date columns (C:G) for the rows that only had Name/Year so far,
# and append two brand new rows (Perl, Rust).
# ---------------------------------------------------------------------------

# Row 6 - Java
$ws.Range("C6").Value = "Джеймс Гослинг"
$ws.Range("D6").Value = "да"
$ws.Range("E6").Value = "да"
$ws.Range("F6").Value = "Java SE 19.0.1"
$ws.Range("G6").Value = "18.10.2022"

# Row 7 - Swift
$ws.Range("C7").Value = "Крис Латтнер"
$ws.Range("D7").Value = "да"
$ws.Range("E7").Value = "да"
$ws.Range("F7").Value = "5.7.2"
$ws.Range("G7").Value = "14.12.2022"

# Row 8 - Groovy
$ws.Range("C8").Value = "Джеймс Стрэкан"
$ws.Range("D8").Value = "да"
$ws.Range("E8").Value = "да"
$ws.Range("F8").Value = "4.0.7"
$ws.Range("G8").Value = "21.12.2022"

# Row 9 - Golang (Go)
$ws.Range("C9").Value = "Роберт Гризмер, Роб Пайк и Кен Томпсон"
$ws.Range("D9").Value = "да"
$ws.Range("E9").Value = "да"
$ws.Range("F9").Value = "1.19.4"
$ws.Range("G9").Value = "06.12.2022"

# Row 10 - Scala
$ws.Range("C10").Value = "Мартин Одерски"
$ws.Range("D10").Value = "да"
$ws.Range("E10").Value = "да"
$ws.Range("F10").Value = "3.1.3"
$ws.Range("G10").Value = "21.06.2022"

# Row 11 - PHP
$ws.Range("C11").Value = "Расмус Лердорф, Энди Гутманс"
$ws.Range("D11").Value = "нет"
$ws.Range("E11").Value = "да"
$ws.Range("F11").Value = "8.2.0"
$ws.Range("G11").Value = "08.12.2022"

# Row 12 - Ruby
$ws.Range("C12").Value = "Юкихиро Мацумото"
$ws.Range("D12").Value = "да"
$ws.Range("E12").Value = "да"
$ws.Range("F12").Value = "3.2.0"
$ws.Range("G12").Value = "25.12.2022"

# Row 13 - Objective-C
$ws.Range("C13").Value = "Бред Кокс"
$ws.Range("D13").Value = "да"
$ws.Range("E13").Value = "нет"
$ws.Range("F13").Value = "2.0"
$ws.Range("G13").Value = "19.10.2022"

# Row 14 - Haskell
$ws.Range("C14").Value = "Леннарт Аугустссон, Уоррен Бертон"
$ws.Range("D14").Value = "нет"
$ws.Range("E14").Value = "нет"
$ws.Range("F14").Value = "Haskell 2010"
$ws.Range("G14").Value = "01.07.2010"

# Row 15 - C++
$ws.Range("C15").Value = "Бьёрн Страуструп"
$ws.Range("D15").Value = "да"
$ws.Range("E15").Value = "да"
$ws.Range("F15").Value = "C++20"
$ws.Range("G15").Value = "01.12.2020"

# Row 16 - JavaScript
$ws.Range("C16").Value = "Брендан Эйх"
$ws.Range("D16").Value = "да"
$ws.Range("E16").Value = "да"
$ws.Range("F16").Value = "ECMAScript 2022"
$ws.Range("G16").Value = "01.06.2022"

# Row 17 - Perl (new row)
$ws.Range("A17").Value = "Perl"
$ws.Range("B17").Value = "1987"
$ws.Range("C17").Value = "Ларри Уолл"
$ws.Range("D17").Value = "нет"
$ws.Range("E17").Value = "да"
$ws.Range("F17").Value = "5.36.0"
$ws.Range("G17").Value = "28.05.2022"

# Row 18 - Rust (new row)
$ws.Range("A18").Value = "Rust"
$ws.Range("B18").Value = "2006"
$ws.Range("C18").Value = "Грэйдон Хор"
$ws.Range("D18").Value = "да"
$ws.Range("E18").Value = "да"
$ws.Range("F18").Value = "1.66.1"
$ws.Range("G18").Value = "10.01.2023"

# ---------------------------------------------------------------------------
# Re-apply the text/wrap formatting to the freshly written cells (Value
# assignment resets a new cell back to the default "General" style) and
# size the rows so the wrapped text is fully visible, matching how Excel
# would auto-fit these rows after the text was entered.
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A6:G18")
$dataRange.NumberFormat = "@"
$dataRange.WrapText = $true

$ws.Rows.Item(6).RowHeight = 28.8
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 57.6
$ws.Rows.Item(10).RowHeight = 28.8
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 28.8
$ws.Rows.Item(14).RowHeight = 57.6
$ws.Rows.Item(15).RowHeight = 28.8
$ws.Rows.Item(16).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Match the final cursor/selection position left by the editing session.
# ---------------------------------------------------------------------------
[void]$ws.Range("F20").Select()
